# Apply the "Saldo" export update:
#   - Row 2 (ANILSON) balance changes 154097.32 -> 108615.2
#   - Row 3 (PRISCILLA / 004224284 / 16030.59) is replaced by DIEGO / 004479965 / 17432.65
#   - Row 4 (CESAR / 004207278 / 9176.22) is replaced by ROSANGELA / 005428871 / 14879.47
#   - A brand-new row is inserted after that for E3 / 004267976 / 11892.73,
#     pushing the remaining untouched rows (MARCELO, GUSTAVO, ...) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the string to be stored as text so values that look numeric
    # (e.g. account numbers with leading zeros) keep their exact digits,
    # then strip the residual number-format/quote-prefix styling so the
    # cell ends up with the same "no explicit style" look as its neighbours.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# 1) Update ANILSON's balance in place.
$ws.Range("C2").Value = 108615.2

# 2) Turn the PRISCILLA row into the DIEGO row.
Set-TextCell $ws.Range("A3") "004479965"
Set-TextCell $ws.Range("B3") "DIEGO"
$ws.Range("C3").Value = 17432.65

# 3) Turn the CESAR row into the ROSANGELA row.
Set-TextCell $ws.Range("A4") "005428871"
Set-TextCell $ws.Range("B4") "ROSANGELA"
$ws.Range("C4").Value = 14879.47

# 4) Insert a brand-new row for E3, pushing MARCELO and everything below
#    down by one row (matches the net +1 row growth seen in the diff).
$ws.Rows(5).Insert()
Set-TextCell $ws.Range("A5") "004267976"
Set-TextCell $ws.Range("B5") "E3"
$ws.Range("C5").Value = 11892.73
